# Apply updated leve-profit figures pulled from the latest market-board snapshot.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1924294
$ws.Range("J17").Value = 1924294
$ws.Range("L17").Value = 5772882
$ws.Range("N17").Value = -5773218

$ws.Range("H29").Value = 1100.8
$ws.Range("I29").Value = 201.6
$ws.Range("J29").Value = 2000
$ws.Range("K29").Value = 604.8
$ws.Range("L29").Value = 6000
$ws.Range("M29").Value = -323.8
$ws.Range("N29").Value = -6562

$ws.Range("H38").Value = 414
$ws.Range("I38").Value = 80
$ws.Range("J38").Value = 1750
$ws.Range("K38").Value = 240
$ws.Range("L38").Value = 5250
$ws.Range("M38").Value = 132
$ws.Range("N38").Value = -5994

$ws.Range("H106").Value = 1820
$ws.Range("I106").Value = 1525
$ws.Range("J106").Value = 3000
$ws.Range("K106").Value = 1525
$ws.Range("L106").Value = 3000
$ws.Range("M106").Value = -894
$ws.Range("N106").Value = -4262

$ws.Range("H132").Value = 30103.756
$ws.Range("I132").Value = 45405.793
$ws.Range("J132").Value = 1853.8462
$ws.Range("K132").Value = 136217.379
$ws.Range("L132").Value = 5561.5386
$ws.Range("M132").Value = -133687.379
$ws.Range("N132").Value = -10621.5386

$ws.Range("H138").Value = 3400.06
$ws.Range("I138").Value = 2307.2778
$ws.Range("J138").Value = 4014.75
$ws.Range("K138").Value = 6921.8334
$ws.Range("L138").Value = 12044.25
$ws.Range("M138").Value = -1781.8334
$ws.Range("N138").Value = -22324.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7243.1465
$ws.Range("I32").Value = 6499.225
$ws.Range("J32").Value = 37000
$ws.Range("K32").Value = 6499.225
$ws.Range("L32").Value = 37000
$ws.Range("M32").Value = -6212.225
$ws.Range("N32").Value = -37574

$ws.Range("H61").Value = 2025.8889
$ws.Range("I61").Value = 2030.3077
$ws.Range("K61").Value = 2030.3077
$ws.Range("M61").Value = -1818.3077

$ws.Range("H74").Value = 1565.6666
$ws.Range("I74").Value = 1156.9584
$ws.Range("J74").Value = 2655.5557
$ws.Range("K74").Value = 1156.9584
$ws.Range("L74").Value = 2655.5557
$ws.Range("M74").Value = -282.9584
$ws.Range("N74").Value = -4403.5557

$ws.Range("H77").Value = 1565.6666
$ws.Range("I77").Value = 1156.9584
$ws.Range("J77").Value = 2655.5557
$ws.Range("K77").Value = 5784.791999999999
$ws.Range("L77").Value = 13277.7785
$ws.Range("M77").Value = -1416.791999999999
$ws.Range("N77").Value = -22013.7785

$ws.Range("H102").Value = 1764.2858
$ws.Range("I102").Value = 1132.5
$ws.Range("K102").Value = 1132.5
$ws.Range("M102").Value = 489.5

$ws.Range("H122").Value = 2105.7
$ws.Range("I122").Value = 1708.3334
$ws.Range("J122").Value = 2701.75
$ws.Range("K122").Value = 5125.0002
$ws.Range("L122").Value = 8105.25
$ws.Range("M122").Value = -2675.0002
$ws.Range("N122").Value = -13005.25

$ws.Range("H132").Value = 1646.4755
$ws.Range("I132").Value = 1326.7142
$ws.Range("J132").Value = 2952.1667
$ws.Range("K132").Value = 3980.1426
$ws.Range("L132").Value = 8856.500100000001
$ws.Range("M132").Value = -1450.1426
$ws.Range("N132").Value = -13916.5001

$ws.Range("H136").Value = 2025.8889
$ws.Range("I136").Value = 2030.3077
$ws.Range("K136").Value = 6090.9231
$ws.Range("M136").Value = -3540.9231

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1129.15
$ws.Range("I99").Value = 1166.9375
$ws.Range("J99").Value = 978
$ws.Range("K99").Value = 1166.9375
$ws.Range("L99").Value = 978
$ws.Range("M99").Value = 331.0625
$ws.Range("N99").Value = -3974

$ws.Range("H105").Value = 2019.3334
$ws.Range("I105").Value = 1986.1538
$ws.Range("J105").Value = 2235
$ws.Range("K105").Value = 1986.1538
$ws.Range("L105").Value = 2235
$ws.Range("M105").Value = -239.1538
$ws.Range("N105").Value = -5729

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H103").Value = 3848.3333
$ws.Range("I103").Value = 45
$ws.Range("J103").Value = 5750
$ws.Range("K103").Value = 135
$ws.Range("L103").Value = 17250
$ws.Range("M103").Value = 744
$ws.Range("N103").Value = -19008

$ws.Range("H131").Value = 772.8484999999999
$ws.Range("I131").Value = 309.23077
$ws.Range("J131").Value = 842.93024
$ws.Range("K131").Value = 927.69231
$ws.Range("L131").Value = 2528.79072
$ws.Range("M131").Value = 4112.30769
$ws.Range("N131").Value = -12608.79072

$ws.Range("H132").Value = 2012
$ws.Range("I132").Value = 998.5
$ws.Range("J132").Value = 2265.375
$ws.Range("K132").Value = 8986.5
$ws.Range("L132").Value = 20388.375
$ws.Range("M132").Value = -6456.5
$ws.Range("N132").Value = -25448.375

$ws.Range("H140").Value = 2324.3547
$ws.Range("I140").Value = 1403.2609
$ws.Range("K140").Value = 4209.7827
$ws.Range("M140").Value = 970.2173000000003

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H49").Value = 0
$ws.Range("J49").Value = 0
$ws.Range("L49").Value = 0
$ws.Range("N49").Value = $null

$ws.Range("H102").Value = 4050890.2
$ws.Range("I102").Value = 4809901
$ws.Range("J102").Value = 2833.3333
$ws.Range("K102").Value = 4809901
$ws.Range("L102").Value = 2833.3333
$ws.Range("M102").Value = -4808279
$ws.Range("N102").Value = -6077.3333

$ws.Range("H122").Value = 3231.4736
$ws.Range("I122").Value = 2171.3333
$ws.Range("K122").Value = 6513.999899999999
$ws.Range("M122").Value = -4063.999899999999

$ws.Range("H132").Value = 2631.1724
$ws.Range("I132").Value = 2439.3845
$ws.Range("J132").Value = 2787
$ws.Range("K132").Value = 7318.1535
$ws.Range("L132").Value = 8361
$ws.Range("M132").Value = -4788.1535
$ws.Range("N132").Value = -13421

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1004.2727
$ws.Range("I22").Value = 1316.6666
$ws.Range("J22").Value = 887.125
$ws.Range("K22").Value = 1316.6666
$ws.Range("L22").Value = 887.125
$ws.Range("M22").Value = -1021.6666
$ws.Range("N22").Value = -1477.125

$ws.Range("H27").Value = 1004.2727
$ws.Range("I27").Value = 1316.6666
$ws.Range("J27").Value = 887.125
$ws.Range("K27").Value = 1316.6666
$ws.Range("L27").Value = 887.125
$ws.Range("M27").Value = -1209.6666
$ws.Range("N27").Value = -1101.125

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2114.4119
$ws.Range("I132").Value = 1781.0646
$ws.Range("J132").Value = 2631.1
$ws.Range("K132").Value = 5343.1938
$ws.Range("L132").Value = 7893.299999999999
$ws.Range("M132").Value = -2813.1938
$ws.Range("N132").Value = -13953.3
